$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 30 (shifts old row 30 -> 31, old row 31 -> 32).
# The new row inherits formatting from the row above it, which already
# matches the table's styling (text columns + blank-date column), so no
# extra style tweaking is required.
$ws.Rows("30").Insert()
$ws.Range("A30").RowHeight = 13.05

# Populate the new row 30 with the new customer record (no invoice date yet)
$ws.Range("A30").Value = "DG PROCESSING"
$ws.Range("B30").Value = "Zigan, Gerald L"
$ws.Range("C30").Value = "030"
$ws.Range("E30").Value = "0008370"

# Column F stays blank for this row, but the cell itself should still be
# materialised (as it is for every other row in the table).
$ws.Range("F30").Font.Name = $ws.Range("F30").Font.Name
